# "Change DB" - add two new root/suffix entries to the Лист3 dictionary sheet
# and repoint an existing row's reference id, per the commit "Added answer.py - Change DB".

$wb = $excel.ActiveWorkbook

# Лист3 is the first worksheet (root/suffix -> id table).
$ws = $wb.Worksheets.Item(1)

# Existing row 7 ("бо") now points to id 14 instead of 5.
$ws.Range("C7").Value = 14

# New row 13: root "щ", suffix "о" (same suffix used in row 11 "лих"/"о"), id 11.
$ws.Range("A13").Value = "щ"
$ws.Range("B13").Value = "о"
$ws.Range("C13").Value = 11

# New row 14: root "чому" (no suffix), id 13.
$ws.Range("A14").Value = "чому"
$ws.Range("C14").Value = 13

# Update the saved selection on this sheet to match the edited area.
$ws.Range("G12").Select() | Out-Null
